# Update transition-probability rows on the "strategy_id-0" sheet and
# move the active selection/scroll position, per the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")
$ws.Activate()

# Row -> new uniform value for columns J:AS (value is repeated across
# every column from J through AS in each of these rows).
$rowValues = @{
    8  = 0.3
    9  = 0.2
    10 = 0.2
    11 = 0.25
    12 = 0.2
    13 = 0.25
    14 = 0.2
    15 = 0.2
    16 = 0.2
    17 = 0.1
    18 = 0.05
    19 = 0.25
    20 = 0.25
    77 = 0
}

foreach ($row in $rowValues.Keys) {
    $value = $rowValues[$row]
    $rng = $ws.Range("J$row" + ":AS$row")
    $rng.Value = $value
}

# Move the saved view state: scroll position + active selection.
[void]$ws.Range("A154").Select()
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H8").Select()
